# Applies the 2026-02-01 16:02:30 scrape refresh to the LP1912 / LP1912-215 / 6203-6173 sheets.
# Values below were derived from the authoritative before/after cell-level diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet 1: LP1912 -----
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Última actualización: 16:02:30"

$ws.Range("A3").Value = "Total filas: 220"

$ws.Range("C15").Value = "215A_EL PATO"

$ws.Range("C16").Value = "225_GOMEZ"

$ws.Range("A35").Value = "07:26:49"
$ws.Range("C35").Value = "16_SANTA ANA"
$ws.Range("D35").Value = 34

$ws.Range("A36").Value = "07:51:40"
$ws.Range("C36").Value = "17_ROMERO"
$ws.Range("D36").Value = 9

$ws.Range("A43").Value = "08:14:55"
$ws.Range("C43").Value = "11_ETCHEVERRY"
$ws.Range("D43").Value = 14

$ws.Range("A44").Value = "06:58:58"
$ws.Range("C44").Value = "15_ABASTO"
$ws.Range("D44").Value = 90

$ws.Range("C117").Value = "17_ROMERO"

$ws.Range("C118").Value = "16_SANTA ANA"

$ws.Range("C137").Value = "15X38_ABASTO"

$ws.Range("C138").Value = "14_ABASTO"

$ws.Range("A151").Value = "13:28:27"
$ws.Range("C151").Value = "215A_EL PATO"
$ws.Range("D151").Value = 5

$ws.Range("A152").Value = "12:43:13"
$ws.Range("C152").Value = "14_ABASTO"
$ws.Range("D152").Value = 50

$ws.Range("C190").Value = "10_OLMOS"

$ws.Range("C192").Value = "15X38_ABASTO"

$ws.Range("A197").Value = "16:02:30"
$ws.Range("B197").Value = "16:02"
$ws.Range("D197").Value = 0

$ws.Range("A198").Value = "16:02:30"
$ws.Range("B198").Value = "16:04"
$ws.Range("D198").Value = 2

$ws.Range("A199").Value = "16:02:30"
$ws.Range("B199").Value = "16:05"
$ws.Range("C199").Value = "16_SANTA ANA"
$ws.Range("D199").Value = 3

$ws.Range("A200").Value = "14:58:38"
$ws.Range("B200").Value = "16:05"
$ws.Range("C200").Value = "14_ABASTO"
$ws.Range("D200").Value = 67

$ws.Range("A201").Value = "15:34:15"
$ws.Range("B201").Value = "16:06"
$ws.Range("C201").Value = "14_ABASTO"
$ws.Range("D201").Value = 32

$ws.Range("A202").Value = "16:02:30"
$ws.Range("B202").Value = "16:14"
$ws.Range("C202").Value = "17_ROMERO"
$ws.Range("D202").Value = 12

$ws.Range("A203").Value = "14:44:54"
$ws.Range("B203").Value = "16:16"
$ws.Range("C203").Value = "10_OLMOS"
$ws.Range("D203").Value = 92

$ws.Range("A204").Value = "16:02:30"
$ws.Range("B204").Value = "16:17"
$ws.Range("C204").Value = "10_OLMOS"
$ws.Range("D204").Value = 15

$ws.Range("B205").Value = "16:18"
$ws.Range("C205").Value = "23_HERNANDEZ"
$ws.Range("D205").Value = 44

$ws.Range("A206").Value = "16:02:30"
$ws.Range("B206").Value = "16:21"
$ws.Range("C206").Value = "23_HERNANDEZ"
$ws.Range("D206").Value = 19

$ws.Range("A207").Value = "16:02:30"
$ws.Range("B207").Value = "16:29"
$ws.Range("C207").Value = "14_ABASTO"
$ws.Range("D207").Value = 27

$ws.Range("B208").Value = "16:30"
$ws.Range("C208").Value = "14_ABASTO"
$ws.Range("D208").Value = 56

$ws.Range("A209").Value = "16:02:30"
$ws.Range("B209").Value = "16:34"
$ws.Range("C209").Value = "83_ALUAR"
$ws.Range("D209").Value = 32

$ws.Range("A210").Value = "14:44:54"
$ws.Range("B210").Value = "16:40"
$ws.Range("C210").Value = "225_GOMEZ"
$ws.Range("D210").Value = 116

$ws.Range("A211").Value = "16:02:30"
$ws.Range("B211").Value = "16:41"
$ws.Range("C211").Value = "225_GOMEZ"
$ws.Range("D211").Value = 39

$ws.Range("B212").Value = "16:46"
$ws.Range("C212").Value = "17_ROMERO"
$ws.Range("D212").Value = 72

$ws.Range("B213").Value = "16:52"
$ws.Range("C213").Value = "16_SANTA ANA"
$ws.Range("D213").Value = 78

$ws.Range("A214").Value = "16:02:30"
$ws.Range("B214").Value = "16:53"
$ws.Range("C214").Value = "11_ETCHEVERRY"
$ws.Range("D214").Value = 51

$ws.Range("B215").Value = "16:54"
$ws.Range("C215").Value = "11_ETCHEVERRY"
$ws.Range("D215").Value = 80

$ws.Range("A216").Value = "16:02:30"
$ws.Range("B216").Value = "16:58"
$ws.Range("C216").Value = "15_ABASTO"
$ws.Range("D216").Value = 56

$ws.Range("A217").Value = "16:02:30"
$ws.Range("B217").Value = "17:07"
$ws.Range("C217").Value = "16_P MOR-SANTA ANA"
$ws.Range("D217").Value = 65
$ws.Range("E217").Value = "LP1912"

$ws.Range("A218").Value = "16:02:30"
$ws.Range("B218").Value = "17:10"
$ws.Range("C218").Value = "215C_EL PATO"
$ws.Range("D218").Value = 68
$ws.Range("E218").Value = "LP1912"

$ws.Range("A219").Value = "16:02:30"
$ws.Range("B219").Value = "17:17"
$ws.Range("C219").Value = "23_HERNANDEZ"
$ws.Range("D219").Value = 75
$ws.Range("E219").Value = "LP1912"

$ws.Range("A220").Value = "16:02:30"
$ws.Range("B220").Value = "17:21"
$ws.Range("C220").Value = "15X38_ABASTO"
$ws.Range("D220").Value = 79
$ws.Range("E220").Value = "LP1912"

$ws.Range("A221").Value = "16:02:30"
$ws.Range("B221").Value = "17:34"
$ws.Range("C221").Value = "17_ROMERO"
$ws.Range("D221").Value = 92
$ws.Range("E221").Value = "LP1912"

$ws.Range("A222").Value = "16:02:30"
$ws.Range("B222").Value = "17:36"
$ws.Range("C222").Value = "27_EL RETIRO"
$ws.Range("D222").Value = 94
$ws.Range("E222").Value = "LP1912"

$ws.Range("A223").Value = "16:02:30"
$ws.Range("B223").Value = "17:38"
$ws.Range("C223").Value = "215B_EL PATO"
$ws.Range("D223").Value = 96
$ws.Range("E223").Value = "LP1912"

$ws.Range("A224").Value = "16:02:30"
$ws.Range("B224").Value = "17:44"
$ws.Range("C224").Value = "215_EL PELIGRO"
$ws.Range("D224").Value = 102
$ws.Range("E224").Value = "LP1912"

$ws.Range("A225").Value = "16:02:30"
$ws.Range("B225").Value = "17:49"
$ws.Range("C225").Value = "10_OLMOS"
$ws.Range("D225").Value = 107
$ws.Range("E225").Value = "LP1912"

# ----- Sheet 2: LP1912-215 -----
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Última actualización: 16:02:30"

$ws.Range("A3").Value = "Total filas: 36"

$ws.Range("A39").Value = "16:02:30"
$ws.Range("D39").Value = 68

$ws.Range("A40").Value = "16:02:30"
$ws.Range("B40").Value = "17:38"
$ws.Range("C40").Value = "215B_EL PATO"
$ws.Range("D40").Value = 96
$ws.Range("E40").Value = "LP1912"

$ws.Range("A41").Value = "16:02:30"
$ws.Range("B41").Value = "17:44"
$ws.Range("C41").Value = "215_EL PELIGRO"
$ws.Range("D41").Value = 102
$ws.Range("E41").Value = "LP1912"

# ----- Sheet 3: 6203-6173 -----
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Última actualización: 16:02:30"

$ws.Range("A37").Value = "16:02:30"
$ws.Range("D37").Value = 0

$ws.Range("A39").Value = "16:02:30"
$ws.Range("D39").Value = 28

$ws.Range("A40").Value = "16:02:30"
$ws.Range("D40").Value = 64

